$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 13350.4
$ws.Range("I34").Value = 4813
$ws.Range("K34").Value = 4813
$ws.Range("M34").Value = -4610
$ws.Range("H36").Value = 13350.4
$ws.Range("I36").Value = 4813
$ws.Range("K36").Value = 4813
$ws.Range("M36").Value = -4098
$ws.Range("H123").Value = 46759.6
$ws.Range("J123").Value = 46759.6
$ws.Range("L123").Value = 46759.6
$ws.Range("N123").Value = -56559.6
$ws.Range("H130").Value = 34990
$ws.Range("J130").Value = 34990
$ws.Range("L130").Value = 34990
$ws.Range("N130").Value = -45030

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 25000
$ws.Range("J55").Value = 25000
$ws.Range("L55").Value = 25000
$ws.Range("N55").Value = -25630
$ws.Range("H80").Value = 38722
$ws.Range("J80").Value = 38722
$ws.Range("L80").Value = 38722
$ws.Range("N80").Value = -40718
$ws.Range("H83").Value = 38722
$ws.Range("J83").Value = 38722
$ws.Range("L83").Value = 116166
$ws.Range("N83").Value = -126150
$ws.Range("H103").Value = 36448
$ws.Range("J103").Value = 36448
$ws.Range("L103").Value = 36448
$ws.Range("N103").Value = -38792
$ws.Range("H109").Value = 27433.334
$ws.Range("J109").Value = 27433.334
$ws.Range("L109").Value = 27433.334
$ws.Range("N109").Value = -30207.334
$ws.Range("H125").Value = 128607820
$ws.Range("J125").Value = 128607820
$ws.Range("L125").Value = 128607820
$ws.Range("N125").Value = -128617660
$ws.Range("H127").Value = 48500
$ws.Range("J127").Value = 48500
$ws.Range("L127").Value = 48500
$ws.Range("H128").Value = 56500
$ws.Range("J128").Value = 56500
$ws.Range("L128").Value = 56500
$ws.Range("H131").Value = 39892.832
$ws.Range("J131").Value = 39892.832
$ws.Range("L131").Value = 39892.832
$ws.Range("N131").Value = -49972.832
$ws.Range("H133").Value = 48000
$ws.Range("J133").Value = 48000
$ws.Range("L133").Value = 48000
$ws.Range("N133").Value = -53060
$ws.Range("H135").Value = 42538.25
$ws.Range("J135").Value = 42538.25
$ws.Range("L135").Value = 42538.25
$ws.Range("N135").Value = -52678.25
$ws.Range("N127").Value = -58420
$ws.Range("N128").Value = -66460

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H112").Value = 28806.9
$ws.Range("J112").Value = 28806.9
$ws.Range("L112").Value = 28806.9
$ws.Range("N112").Value = -31760.9
$ws.Range("H130").Value = 39115
$ws.Range("J130").Value = 39115
$ws.Range("L130").Value = 39115
$ws.Range("N130").Value = -49155
$ws.Range("H135").Value = 53615.24
$ws.Range("J135").Value = 53615.24
$ws.Range("L135").Value = 53615.24
$ws.Range("N135").Value = -63755.24
$ws.Range("H137").Value = 39829.69
$ws.Range("J137").Value = 39829.69
$ws.Range("L137").Value = 39829.69
$ws.Range("N137").Value = -50029.69

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 48999.668
$ws.Range("J20").Value = 48999.668
$ws.Range("L20").Value = 48999.668
$ws.Range("N20").Value = -49471.668
$ws.Range("H30").Value = 48999.668
$ws.Range("J30").Value = 48999.668
$ws.Range("L30").Value = 48999.668
$ws.Range("N30").Value = -49181.668
$ws.Range("H31").Value = 2070.2253
$ws.Range("I31").Value = 1457.3864
$ws.Range("J31").Value = 3068.926
$ws.Range("K31").Value = 1457.3864
$ws.Range("L31").Value = 3068.926
$ws.Range("M31").Value = -1162.3864
$ws.Range("N31").Value = -3658.926
$ws.Range("H34").Value = 2070.2253
$ws.Range("I34").Value = 1457.3864
$ws.Range("J34").Value = 3068.926
$ws.Range("K34").Value = 1457.3864
$ws.Range("L34").Value = 3068.926
$ws.Range("M34").Value = -1255.3864
$ws.Range("N34").Value = -3472.926
$ws.Range("H109").Value = 11000
$ws.Range("J109").Value = 11000
$ws.Range("L109").Value = 11000
$ws.Range("N109").Value = -13080
$ws.Range("H124").Value = 27457.5
$ws.Range("J124").Value = 27457.5
$ws.Range("L124").Value = 27457.5
$ws.Range("N124").Value = -32367.5
$ws.Range("H128").Value = 48999.668
$ws.Range("J128").Value = 48999.668
$ws.Range("L128").Value = 48999.668
$ws.Range("N128").Value = -58959.668
$ws.Range("H135").Value = 55824
$ws.Range("J135").Value = 55824
$ws.Range("L135").Value = 55824
$ws.Range("N135").Value = -65964

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 908.08
$ws.Range("J131").Value = 926.5851
$ws.Range("L131").Value = 2779.7553
$ws.Range("N131").Value = -12859.7553

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 26687
$ws.Range("J57").Value = 26687
$ws.Range("L57").Value = 26687
$ws.Range("N57").Value = -28327
$ws.Range("H93").Value = 9818.75
$ws.Range("J93").Value = 9818.75
$ws.Range("L93").Value = 9818.75
$ws.Range("N93").Value = -13562.75
$ws.Range("H94").Value = 31700
$ws.Range("J94").Value = 31700
$ws.Range("L94").Value = 31700
$ws.Range("N94").Value = -33052
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("H102").Value = 3707263.8
$ws.Range("I102").Value = 7408540.5
$ws.Range("J102").Value = 5986.6665
$ws.Range("K102").Value = 7408540.5
$ws.Range("L102").Value = 5986.6665
$ws.Range("M102").Value = -7406918.5
$ws.Range("N102").Value = -9230.666499999999
$ws.Range("H105").Value = 31853
$ws.Range("J105").Value = 31853
$ws.Range("L105").Value = 31853
$ws.Range("N105").Value = -38841
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("H127").Value = 48475.332
$ws.Range("J127").Value = 48475.332
$ws.Range("L127").Value = 48475.332
$ws.Range("N127").Value = -58395.332
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("N124").ClearContents()
$ws.Range("N128").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1903.75
$ws.Range("J40").Value = 1900
$ws.Range("L40").Value = 1900
$ws.Range("N40").Value = -2172
$ws.Range("H63").Value = 27336.666
$ws.Range("J63").Value = 27336.666
$ws.Range("L63").Value = 27336.666
$ws.Range("N63").Value = -28834.666
$ws.Range("H66").Value = 27336.666
$ws.Range("J66").Value = 27336.666
$ws.Range("L66").Value = 82009.99800000001
$ws.Range("N66").Value = -89497.99800000001
$ws.Range("H122").Value = 9263135
$ws.Range("I122").Value = 33335494
$ws.Range("K122").Value = 100006482
$ws.Range("M122").Value = -100004032
$ws.Range("H127").Value = 56428.332
$ws.Range("J127").Value = 56428.332
$ws.Range("L127").Value = 56428.332
$ws.Range("N127").Value = -66348.33199999999
$ws.Range("H128").Value = 56250
$ws.Range("J128").Value = 56250
$ws.Range("L128").Value = 56250
$ws.Range("N128").Value = -66210
$ws.Range("H130").Value = 56395
$ws.Range("J130").Value = 56395
$ws.Range("L130").Value = 56395
$ws.Range("N130").Value = -66435
$ws.Range("H133").Value = 46326.125
$ws.Range("J133").Value = 46326.125
$ws.Range("L133").Value = 46326.125
$ws.Range("N133").Value = -51386.125

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 12482.25
$ws.Range("I122").Value = 14540.9375
$ws.Range("J122").Value = 4247.5
$ws.Range("K122").Value = 43622.8125
$ws.Range("L122").Value = 12742.5
$ws.Range("M122").Value = -41172.8125
$ws.Range("N122").Value = -17642.5
$ws.Range("H127").Value = 19879
$ws.Range("J127").Value = 19879
$ws.Range("L127").Value = 19879
$ws.Range("N127").Value = -29799
$ws.Range("H128").Value = 48699.855
$ws.Range("J128").Value = 48699.855
$ws.Range("L128").Value = 48699.855
$ws.Range("N128").Value = -58659.855
